$wb = $excel.ActiveWorkbook

# Old / new identifiers for this handoff generation.
$oldGuid = "d00873f9-8c3a-4cc6-9427-1052b4a6efe2"
$newGuid = "637ea68b-9530-496a-8634-572befa58fe0"

$newMdName     = "$newGuid.md"
$newPathName   = "e2e\$newGuid.md"
$newZhCnXlf    = "$newGuid.b4dccdb3cbd5d4f2873307003f8f6c4628faa669.zh-cn.xlf"
$newDeDeXlf    = "$newGuid.b4dccdb3cbd5d4f2873307003f8f6c4628faa669.de-de.xlf"

$newGenerateDate = "2016-09-04 21:04:43"
$newZhCnHandoffDate = "2016-09-04 21:04:36"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newPathName
$wsOverview.Range("G2").Value = $newGenerateDate
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $newPathName
}

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newZhCnHandoffDate
foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = $newMdName
}

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newGenerateDate
foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = $newMdName
}
